$d = $word.ActiveDocument

$pairs = @(
    @("961×7=6727", "938×9=8442"),
    @("912×2=1824", "744×3=2232"),
    @("533×3=1599", "230×4=920"),
    @("881×3=2643", "113×3=339"),
    @("905×9=8145", "732×3=2196"),
    @("844×2=1688", "678×4=2712"),
    @("680×5=3400", "432×2=864"),
    @("265×3=795", "984×7=6888"),
    @("136×4=544", "288×2=576"),
    @("251×6=1506", "783×8=6264"),
    @("861×2=1722", "739×4=2956"),
    @("901×8=7208", "635×3=1905"),
    @("675×2=1350", "281×4=1124"),
    @("509×8=4072", "646×7=4522"),
    @("301×5=1505", "825×4=3300"),
    @("999×3=2997", "777×2=1554"),
    @("622×8=4976", "235×2=470"),
    @("828×4=3312", "629×3=1887"),
    @("197×8=1576", "918×7=6426"),
    @("534×8=4272", "475×5=2375"),
    @("963×4=3852", "376×9=3384"),
    @("419×9=3771", "242×2=484"),
    @("109×9=981", "630×8=5040"),
    @("622×6=3732", "757×8=6056"),
    @("330×7=2310", "297×2=594")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
